# Insert a new data row at row 385 (shifting the existing rows 385-480 down
# to 386-481), then populate the new row 385 with its values. This mirrors
# the diff, which effectively prepends one new price record ahead of the
# existing row that used to be at 385 (and every following row), growing the
# sheet from A1:R480 to A1:R481.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 385:480 down to 386:481, leaving a blank row 385 to fill in.
$ws.Rows(385).Insert()

# Populate the new row 385 with the new record's data.
$ws.Cells.Item(385, 1).Value  = 5
$ws.Cells.Item(385, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(385, 3).Value  = "Maule"
$ws.Cells.Item(385, 4).Value  = 44932
$ws.Cells.Item(385, 5).Value  = 7
$ws.Cells.Item(385, 6).Value  = 100112032
$ws.Cells.Item(385, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(385, 8).Value  = "Sin especificar"
$ws.Cells.Item(385, 9).Value  = "Primera"
$ws.Cells.Item(385, 10).Value = 500
$ws.Cells.Item(385, 11).Value = 4000
$ws.Cells.Item(385, 12).Value = 4000
$ws.Cells.Item(385, 13).Value = 4000
$ws.Cells.Item(385, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(385, 15).Value = "Región del Maule"
$ws.Cells.Item(385, 16).Value = 80
$ws.Cells.Item(385, 17).Value = 50
$ws.Cells.Item(385, 18).Value = "Hortaliza"
